$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 7229
$ws.Range("I3").Value = 7448
$ws.Range("G4").Value = 1449
$ws.Range("I4").Value = 1711
$ws.Range("I5").Value = 704
$ws.Range("I6").Value = 8906
$ws.Range("G7").Value = 24674
$ws.Range("I7").Value = 25998

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I2").Value = 89
$ws.Range("I6").Value = 124
$ws.Range("I7").Value = 305

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I3").Value = 261
$ws.Range("I6").Value = 244
$ws.Range("I7").Value = 804

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 234
$ws.Range("I6").Value = 292
$ws.Range("I7").Value = 975

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 194
$ws.Range("I3").Value = 183
$ws.Range("I7").Value = 609

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I7").Value = 816
$ws.Range("I8").Value = 1537
$ws.Range("I11").Value = 395
$ws.Range("I12").Value = 65
$ws.Range("I17").Value = 40
$ws.Range("I18").Value = 205
$ws.Range("I19").Value = 726
$ws.Range("I23").Value = 251
$ws.Range("I25").Value = 133
$ws.Range("I27").Value = 226
$ws.Range("I29").Value = 1547
$ws.Range("I33").Value = 1139
$ws.Range("I36").Value = 357
$ws.Range("I37").Value = 804
$ws.Range("I42").Value = 1000
$ws.Range("I47").Value = 189
$ws.Range("I48").Value = 330
$ws.Range("I50").Value = 137
$ws.Range("I52").Value = 589
$ws.Range("I54").Value = 506
$ws.Range("I60").Value = 152
$ws.Range("G63").Value = 210
$ws.Range("I63").Value = 85
$ws.Range("I65").Value = 609
$ws.Range("I67").Value = 975
$ws.Range("I70").Value = 47
$ws.Range("I73").Value = 235
$ws.Range("I78").Value = 345
$ws.Range("I79").Value = 745
$ws.Range("I83").Value = 567
$ws.Range("I85").Value = 1154
$ws.Range("I86").Value = 170
$ws.Range("I91").Value = 277
$ws.Range("I93").Value = 149
$ws.Range("I94").Value = 259
$ws.Range("I96").Value = 305
$ws.Range("I98").Value = 188
$ws.Range("G101").Value = 24674
$ws.Range("I101").Value = 25998

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 189
$ws.Range("I3").Value = 205
$ws.Range("I6").Value = 124
$ws.Range("I7").Value = 567

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I4").Value = 17
$ws.Range("I6").Value = 86

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 414
$ws.Range("I7").Value = 1139

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I6").Value = 246
$ws.Range("I7").Value = 506

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value = 529
$ws.Range("I6").Value = 429
$ws.Range("I7").Value = 1547

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value = 212
$ws.Range("I7").Value = 726

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I6").Value = 170
$ws.Range("I7").Value = 330

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 437
$ws.Range("I5").Value = 36
$ws.Range("I6").Value = 302
$ws.Range("I7").Value = 1154

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I6").Value = 420
$ws.Range("I7").Value = 1000

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I6").Value = 127
$ws.Range("I7").Value = 345

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 70
$ws.Range("I7").Value = 251

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I2").Value = 86
$ws.Range("I3").Value = 99
$ws.Range("I7").Value = 277

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I3").Value = 243
$ws.Range("I6").Value = 213
$ws.Range("I7").Value = 745

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I6").Value = 97
$ws.Range("I7").Value = 205

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I3").Value = 119
$ws.Range("I7").Value = 357

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I3").Value = 37
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 145
$ws.Range("I3").Value = 192
$ws.Range("I6").Value = 193
$ws.Range("I7").Value = 589

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I2").Value = 53
$ws.Range("I3").Value = 41
$ws.Range("I7").Value = 259

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I6").Value = 38
$ws.Range("I7").Value = 133

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I2").Value = 47
$ws.Range("I7").Value = 189

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I3").Value = 16
$ws.Range("I7").Value = 188

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I3").Value = 30
$ws.Range("I7").Value = 137

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 155
$ws.Range("I6").Value = 111
$ws.Range("I7").Value = 395

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I6").Value = 63
$ws.Range("I7").Value = 235

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("I2").Value = 15
$ws.Range("I7").Value = 47

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 454
$ws.Range("I3").Value = 444
$ws.Range("I6").Value = 497
$ws.Range("I7").Value = 1537

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I3").Value = 43
$ws.Range("I7").Value = 226

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 170

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I3").Value = 42
$ws.Range("I7").Value = 152

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I6").Value = 227
$ws.Range("I7").Value = 816

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("I6").Value = 35
$ws.Range("I7").Value = 65
